$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'45.286.20"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -1.15%  "
$c = $ws.Range("D3")
$c.Value = "'2.369.17"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -2.34%  "
$ws.Range("E4").Value = "  -0.05%  "
$c = $ws.Range("D5")
$c.Value = "'318.67"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.36%  "
$c = $ws.Range("D6")
$c.Value = "'108.44"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -5.63%  "
$c = $ws.Range("D7")
$c.Value = "'0.635"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("E9").Value = "  -3.01%  "
$c = $ws.Range("D10")
$c.Value = "'40.99"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -5.09%  "
$c = $ws.Range("D11")
$c.Value = "'0.0920"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -2.53%  "
$ws.Range("E12").Value = "  -3.25%  "
$ws.Range("E13").Value = "  +0.01%  "
$c = $ws.Range("D14")
$c.Value = "'0.981"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -4.41%  "
$c = $ws.Range("D15")
$c.Value = "'2.730.11"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -1.77%  "
$c = $ws.Range("D16")
$c.Value = "'15.44"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -3.74%  "
$c = $ws.Range("D17")
$c.Value = "'2.365.18"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -2.25%  "
$c = $ws.Range("D18")
$c.Value = "'45.243.11"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -1.23%  "
$c = $ws.Range("D19")
$c.Value = "'15.52"
$c.Style = "Normal"
$c = $ws.Range("D20")
$c.Value = "'7.29"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -2.55%  "
$c = $ws.Range("D22")
$c.Value = "'3.61"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +1.41%  "
$c = $ws.Range("D23")
$c.Value = "'73.24"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -2.54%  "
$c = $ws.Range("D24")
$c.Value = "'264.73"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -1.85%  "
$ws.Range("E25").Value = "  -1.45%  "
$ws.Range("E26").Value = "  +0.01%  "
$c = $ws.Range("D27")
$c.Value = "'11.24"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -1.14%  "
$c = $ws.Range("D28")
$c.Value = "'7.50"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -1.78%  "
$ws.Range("E29").Value = "  -1.68%  "
$c = $ws.Range("D30")
$c.Value = "'22.38"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -3.12%  "
$c = $ws.Range("D31")
$c.Value = "'0.0948"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -2.18%  "
$ws.Range("E32").Value = "  -4.56%  "
$c = $ws.Range("D33")
$c.Value = "'168.89"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -2.18%  "
$ws.Range("E34").Value = "  -3.91%  "
$c = $ws.Range("D35")
$c.Value = "'0.132"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.61%  "
$ws.Range("E36").Value = "  -4.47%  "
$ws.Range("E37").Value = "  -5.55%  "
$c = $ws.Range("D38")
$c.Value = "'3.08"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -1.38%  "
$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Range("D39")
$c.Value = "'4.04"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -2.38%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c = $ws.Range("D40")
$c.Value = "'1.92"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +6.88%  "
$ws.Range("E41").Value = "  -3.70%  "
$c = $ws.Range("D42")
$c.Value = "'99.18"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -2.97%  "
$c = $ws.Range("D43")
$c.Value = "'70.22"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -2.62%  "
$ws.Range("B44").Value = "Celestia"
$ws.Range("C44").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$c = $ws.Range("D44")
$c.Value = "'12.97"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -3.08%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$c = $ws.Range("D45")
$c.Value = "'1.868.35"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +12.80%  "
$ws.Range("E46").Value = "  -5.34%  "
$ws.Range("E47").Value = "  -0.02%  "
$c = $ws.Range("D48")
$c.Value = "'5.97"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +2.09%  "
$c = $ws.Range("D49")
$c.Value = "'84.55"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +6.89%  "
$c = $ws.Range("D50")
$c.Value = "'112.21"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -5.33%  "
$c = $ws.Range("D51")
$c.Value = "'9.18"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -2.50%  "
